$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 14:53"

# Row 4
$ws.Range("B4").Value = 5747544
$ws.Range("C4").Value = 1272
$ws.Range("D4").Value = 3096317
$ws.Range("E4").Value = 2473771
$ws.Range("G4").Value = 32
$ws.Range("H4").Value = 177456

# Row 6
$ws.Range("B6").Value = 2915015
$ws.Range("C6").Value = 10686
$ws.Range("D6").Value = 2165682
$ws.Range("E6").Value = 694267
$ws.Range("G6").Value = 91
$ws.Range("H6").Value = 55066

# Row 17
$ws.Range("B17").Value = 305186
$ws.Range("C17").Value = 1213
$ws.Range("D17").Value = 277067
$ws.Range("E17").Value = 24539
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 3580

# Row 36
$ws.Range("B36").Value = 86068
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 5810

# Row 46
$ws.Range("B46").Value = 65589
$ws.Range("C46").Value = 535
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 6195

# Row 80
$ws.Range("B80").Value = 16127
$ws.Range("C80").Value = 71
$ws.Range("D80").Value = 13944
$ws.Range("E80").Value = 1562

# Row 102
$ws.Range("B102").Value = 7594
$ws.Range("C102").Value = 265
$ws.Range("D102").Value = 5584
$ws.Range("E102").Value = 1841
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 169

# Row 141
$ws.Range("B141").Value = 1848
$ws.Range("C141").Value = 98
$ws.Range("D141").Value = 1199
$ws.Range("E141").Value = 630

# Row 184
$ws.Range("B184").Value = 231
$ws.Range("C184").Value = 2
$ws.Range("E184").Value = 31
